$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the Price column (D) whose new value looks like a plain
# number need to be forced to Text format first, otherwise Excel
# would auto-convert the assigned string into a numeric value -
# the source workbook stores every Price cell as text.
$textCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D14", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "65.919.58"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "3.525.36"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "603.09"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").Value = "143.21"
$ws.Range("E6").Value = "  -3.48%  "
$ws.Range("D7").Value = "3.524.44"
$ws.Range("E7").Value = "  -1.98%  "
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").Value = "0.509"
$ws.Range("E9").Value = "  +4.07%  "
$ws.Range("D10").Value = "7.83"
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("D11").Value = "0.130"
$ws.Range("E11").Value = "  -4.97%  "
$ws.Range("D12").Value = "0.404"
$ws.Range("E12").Value = "  -2.72%  "
$ws.Range("D13").Value = "4.119.16"
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "0.0000193"
$ws.Range("E14").Value = "  -7.77%  "
$ws.Range("D15").Value = "28.23"
$ws.Range("E15").Value = "  -6.04%  "
$ws.Range("D16").Value = "3.499.35"
$ws.Range("E16").Value = "  -3.88%  "
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").Value = "65.797.86"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").Value = "10.88"
$ws.Range("E19").Value = "  -5.03%  "
$ws.Range("D20").Value = "6.17"
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("D21").Value = "14.48"
$ws.Range("E21").Value = "  -4.06%  "
$ws.Range("D22").Value = "418.37"
$ws.Range("E22").Value = "  -2.92%  "
$ws.Range("D23").Value = "0.593"
$ws.Range("E23").Value = "  -4.59%  "
$ws.Range("D24").Value = "76.62"
$ws.Range("E24").Value = "  -3.24%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "3.649.51"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "0.0000113"
$ws.Range("E27").Value = "  -7.46%  "
$ws.Range("D28").Value = "2.45"
$ws.Range("E28").Value = "  -2.49%  "
$ws.Range("D29").Value = "7.79"
$ws.Range("E29").Value = "  -5.71%  "
$ws.Range("D30").Value = "8.86"
$ws.Range("E30").Value = "  -5.17%  "
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "3.517.74"
$ws.Range("E32").Value = "  -2.09%  "
$ws.Range("D33").Value = "0.154"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "24.16"
$ws.Range("E34").Value = "  -5.60%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "1.35"
$ws.Range("E36").Value = "  -7.58%  "
$ws.Range("D37").Value = "7.48"
$ws.Range("E37").Value = "  -4.82%  "
$ws.Range("D38").Value = "174.37"
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("D39").Value = "1.60"
$ws.Range("E39").Value = "  -6.95%  "
$ws.Range("D40").Value = "5.19"
$ws.Range("E40").Value = "  -7.72%  "
$ws.Range("D41").Value = "0.0811"
$ws.Range("E41").Value = "  -5.33%  "
$ws.Range("D42").Value = "0.855"
$ws.Range("E42").Value = "  -4.75%  "
$ws.Range("D43").Value = "4.93"
$ws.Range("E43").Value = "  -6.01%  "
$ws.Range("D44").Value = "45.45"
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "1.76"
$ws.Range("E45").Value = "  -7.69%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "0.998"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "2.35"
$ws.Range("E47").Value = "  -9.02%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "23.20"
$ws.Range("E48").Value = "  -3.74%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "1.12"
$ws.Range("E49").Value = "  -6.72%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "6.98"
$ws.Range("E50").Value = "  -3.31%  "
$ws.Range("D51").Value = "0.901"
$ws.Range("E51").Value = "  -5.59%  "
